$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.004.13'
$ws.Range("E2").Value = '  +0.15%  '

$ws.Range("D3").Value = '3.912.00'
$ws.Range("E3").Value = '  +5.47%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.46%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '612.38'
$ws.Range("E5").Value = '  -0.95%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.76'
$ws.Range("E6").Value = '  -1.81%  '

$ws.Range("D7").Value = '3.912.65'
$ws.Range("E7").Value = '  +5.39%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.529'
$ws.Range("E9").Value = '  -0.44%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.168'
$ws.Range("E10").Value = '  +1.55%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.49'
$ws.Range("E11").Value = '  +3.26%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.484'
$ws.Range("E12").Value = '  +0.31%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.18'
$ws.Range("E13").Value = '  -0.48%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000256'
$ws.Range("E14").Value = '  +0.01%  '

$ws.Range("D15").Value = '4.524.76'
$ws.Range("E15").Value = '  +4.08%  '

$ws.Range("D16").Value = '3.867.87'
$ws.Range("E16").Value = '  +3.76%  '

$ws.Range("D17").Value = '70.030.88'
$ws.Range("E17").Value = '  +0.09%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.51'
$ws.Range("E18").Value = '  -0.95%  '

$ws.Range("E19").Value = '  -2.92%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.72'
$ws.Range("E20").Value = '  +0.71%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '507.65'
$ws.Range("E21").Value = '  +0.89%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.77'
$ws.Range("E22").Value = '  +5.73%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.745'
$ws.Range("E23").Value = '  +3.89%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.48'
$ws.Range("E24").Value = '  -2.04%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.64'
$ws.Range("E25").Value = '  +0.42%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000145'
$ws.Range("E26").Value = '  +7.72%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.68'
$ws.Range("E27").Value = '  -2.59%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.59'
$ws.Range("E28").Value = '  -6.47%  '

$ws.Range("E29").Value = '  +0.25%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.58'
$ws.Range("E30").Value = '  +4.46%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.00'
$ws.Range("E31").Value = '  +2.87%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '33.15'
$ws.Range("E32").Value = '  +7.09%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.98'
$ws.Range("E33").Value = '  +1.18%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.114'
$ws.Range("E34").Value = '  -0.33%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.996'
$ws.Range("E35").Value = '  -0.56%  '

$ws.Range("E36").Value = '  -0.09%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.16'
$ws.Range("E37").Value = '  +0.64%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.142'
$ws.Range("E38").Value = '  +3.11%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '480.00'
$ws.Range("E39").Value = '  +11.77%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.338'
$ws.Range("E40").Value = '  +0.37%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.07'
$ws.Range("E41").Value = '  +0.21%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '49.81'
$ws.Range("E42").Value = '  -0.38%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.98'
$ws.Range("E43").Value = '  +3.26%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '43.50'
$ws.Range("E44").Value = '  -3.85%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.57'
$ws.Range("E45").Value = '  -0.80%  '

$ws.Range("D46").Value = '2.944.30'
$ws.Range("E46").Value = '  -1.18%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0364'
$ws.Range("E47").Value = '  +0.80%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '139.76'
$ws.Range("E48").Value = '  +2.33%  '

$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '27.28'
$ws.Range("E49").Value = '  -0.35%  '

$ws.Range("B50").Value = 'USDe'
$ws.Range("C50").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.00'
$ws.Range("E50").Value = '  +0.03%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.44'
$ws.Range("E51").Value = '  -1.83%  '
